$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.55

# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("J3").Value = 1.13
$ws.Range("K3").Value = 6
$ws.Range("N3").Value = 3.1
$ws.Range("O3").Value = 1.36

# Row 4
$ws.Range("J4").Value = 1.03
$ws.Range("K4").Value = 9
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 2.35
$ws.Range("T4").Value = 14.5
$ws.Range("U4").Value = 22
$ws.Range("W4").Value = 50
$ws.Range("X4").Value = 26
$ws.Range("Y4").Value = 26
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 7.5
$ws.Range("AB4").Value = 11.75
$ws.Range("AE4").Value = 9.75
$ws.Range("AF4").Value = 11
$ws.Range("AH4").Value = 18
$ws.Range("AI4").Value = 14
$ws.Range("AJ4").Value = 20

# Row 6
$ws.Range("G6").Value = 1.36
$ws.Range("N6").Value = 1.93
$ws.Range("O6").Value = 1.93

# Row 10
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 8
$ws.Range("P10").Value = 1.5
$ws.Range("Q10").Value = 2.5
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 1.75
$ws.Range("T10").Value = 6.5
$ws.Range("Z10").Value = 8
$ws.Range("AD10").Value = 451
$ws.Range("AE10").Value = 8
$ws.Range("AF10").Value = 15
$ws.Range("AI10").Value = 29

# Row 11
$ws.Range("L11").Value = 1.4
$ws.Range("M11").Value = 2.75
$ws.Range("N11").Value = 2.22
$ws.Range("O11").Value = 1.59
$ws.Range("P11").Value = 1.49
$ws.Range("Q11").Value = 2.45
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.73
$ws.Range("T11").Value = 4.9
$ws.Range("U11").Value = 6.8
$ws.Range("V11").Value = 7.2
$ws.Range("W11").Value = 14
$ws.Range("X11").Value = 15
$ws.Range("Y11").Value = 30
$ws.Range("Z11").Value = 6.2
$ws.Range("AA11").Value = 4.9
$ws.Range("AB11").Value = 14
$ws.Range("AC11").Value = 90
$ws.Range("AD11").Value = 101
$ws.Range("AE11").Value = 7
$ws.Range("AF11").Value = 14
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 45
$ws.Range("AI11").Value = 30
$ws.Range("AJ11").Value = 45

# Row 12
$ws.Range("L12").Value = 1.4
$ws.Range("M12").Value = 2.75
$ws.Range("N12").Value = 2.22
$ws.Range("O12").Value = 1.59
$ws.Range("P12").Value = 1.5
$ws.Range("Q12").Value = 2.45
$ws.Range("R12").Value = 1.93
$ws.Range("S12").Value = 1.78
$ws.Range("T12").Value = 5.6
$ws.Range("U12").Value = 8.800000000000001
$ws.Range("V12").Value = 8
$ws.Range("W12").Value = 21
$ws.Range("X12").Value = 19
$ws.Range("Y12").Value = 35
$ws.Range("Z12").Value = 6.2
$ws.Range("AA12").Value = 4.7
$ws.Range("AB12").Value = 13
$ws.Range("AC12").Value = 80
$ws.Range("AD12").Value = 101
$ws.Range("AE12").Value = 6
$ws.Range("AF12").Value = 10
$ws.Range("AG12").Value = 8.6
$ws.Range("AH12").Value = 27
$ws.Range("AI12").Value = 22
$ws.Range("AJ12").Value = 35

# Row 13
$ws.Range("L13").Value = 1.29
$ws.Range("M13").Value = 3.3
$ws.Range("N13").Value = 1.89
$ws.Range("O13").Value = 1.81
$ws.Range("R13").Value = 2.34
$ws.Range("S13").Value = 1.53
$ws.Range("T13").Value = 4.5
$ws.Range("U13").Value = 4.4
$ws.Range("V13").Value = 7.2
$ws.Range("W13").Value = 6.2
$ws.Range("X13").Value = 10
$ws.Range("Y13").Value = 35
$ws.Range("Z13").Value = 8.199999999999999
$ws.Range("AA13").Value = 7.2
$ws.Range("AB13").Value = 23
$ws.Range("AC13").Value = 101
$ws.Range("AD13").Value = 101
$ws.Range("AE13").Value = 15
$ws.Range("AF13").Value = 50
$ws.Range("AG13").Value = 24
$ws.Range("AH13").Value = 101
$ws.Range("AI13").Value = 101
$ws.Range("AJ13").Value = 101

# Row 14
$ws.Range("G14").Value = 2.8
$ws.Range("I14").Value = 2.5

# Row 15
$ws.Range("G15").Value = 2.88
$ws.Range("I15").Value = 2.63

# Row 16
$ws.Range("G16").Value = 2.6
$ws.Range("I16").Value = 2.88
$ws.Range("J16").Value = 1.11
$ws.Range("K16").Value = 6.5
$ws.Range("L16").Value = 1.53
$ws.Range("M16").Value = 2.38
$ws.Range("N16").Value = 2.7
$ws.Range("O16").Value = 1.44

# Row 17
$ws.Range("G17").Value = 1.65

# Row 18
$ws.Range("G18").Value = 2.8
$ws.Range("I18").Value = 2.6

# Row 19
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.7

# Row 20
$ws.Range("G20").Value = 2.3
$ws.Range("L20").Value = 1.36
$ws.Range("M20").Value = 3

# Row 21
$ws.Range("G21").Value = 2.88
$ws.Range("I21").Value = 2.63

# Row 22
$ws.Range("G22").Value = 2.75
$ws.Range("H22").Value = 2.9
$ws.Range("I22").Value = 2.7
$ws.Range("R22").Value = 2.1
$ws.Range("S22").Value = 1.67
$ws.Range("Z22").Value = 6.5
$ws.Range("AF22").Value = 12
$ws.Range("AH22").Value = 29

# Row 23
$ws.Range("G23").Value = 2.88
$ws.Range("I23").Value = 2.25
$ws.Range("L23").Value = 1.2
$ws.Range("M23").Value = 4.33
$ws.Range("N23").Value = 1.67
$ws.Range("O23").Value = 2.15

# Row 24
$ws.Range("I24").Value = 2.1

# Row 25
$ws.Range("I25").Value = 2.3
$ws.Range("L25").Value = 1.44
$ws.Range("M25").Value = 2.63
$ws.Range("N25").Value = 2.4
$ws.Range("O25").Value = 1.53
$ws.Range("T25").Value = 8
$ws.Range("AC25").Value = 67
$ws.Range("AD25").Value = 1250

# Row 26
$ws.Range("I26").Value = 2.05

# Row 27
$ws.Range("G27").Value = 2.7
$ws.Range("I27").Value = 2.4

# Row 28
$ws.Range("G28").Value = 2.63
$ws.Range("I28").Value = 2.63
$ws.Range("J28").Value = 1.1
$ws.Range("K28").Value = 7
$ws.Range("L28").Value = 1.5
$ws.Range("M28").Value = 2.5
$ws.Range("N28").Value = 2.5
$ws.Range("O28").Value = 1.5

# Row 29
$ws.Range("G29").Value = 2.5
$ws.Range("I29").Value = 2.7
$ws.Range("N29").Value = 2.03
$ws.Range("O29").Value = 1.83

# Row 30
$ws.Range("G30").Value = 2.8
$ws.Range("I30").Value = 2.45
$ws.Range("N30").Value = 2.03
$ws.Range("O30").Value = 1.83

# Row 31
$ws.Range("G31").Value = 2.1

# Row 33
$ws.Range("L33").Value = 1.22
$ws.Range("M33").Value = 4
$ws.Range("N33").Value = 1.8
$ws.Range("O33").Value = 2

# Row 39
$ws.Range("G39").Value = 1.29

# Row 40
$ws.Range("G40").Value = 1.95
$ws.Range("I40").Value = 3.7
$ws.Range("N40").Value = 1.98
$ws.Range("O40").Value = 1.88
$ws.Range("P40").Value = 1.4
$ws.Range("Q40").Value = 2.75
$ws.Range("R40").Value = 1.8
$ws.Range("S40").Value = 1.95
$ws.Range("U40").Value = 9.5
$ws.Range("V40").Value = 8.5
$ws.Range("W40").Value = 17
$ws.Range("AD40").Value = 251
$ws.Range("AG40").Value = 13
$ws.Range("AJ40").Value = 41

Write-Output "Applied all changes"